$wb = $excel.ActiveWorkbook

$wsCreate = $wb.Worksheets.Item("Create")
$wsEdit   = $wb.Worksheets.Item("Edit")
$wsDelete = $wb.Worksheets.Item("Delete")

# --- Row 2 data (identical across Create / Edit / Delete) ---
# Processed on "Create" first so the new shared strings land at the same
# indices the target workbook expects (29/30/31).
$wsCreate.Range("A2").Value = "'2265"
$wsCreate.Range("B2").Value = "Mapping"
$wsCreate.Range("C2").Value = "'6756"

$wsEdit.Range("A2").Value = "'2265"
$wsEdit.Range("B2").Value = "Mapping"
$wsEdit.Range("C2").Value = "'6756"

$wsDelete.Range("A2").Value = "'2265"
$wsDelete.Range("B2").Value = "Mapping"
$wsDelete.Range("C2").Value = "'6756"

# --- Row 3 data ---
# "Edit" and "Delete" share identical new values, so do them first (claims
# the next shared-string slots), then "Create" last with its own distinct
# values (claims the final slots) -- this reproduces the exact shared
# string ordering of the target workbook.
$wsEdit.Range("A3").Value = "'9809"
$wsEdit.Range("B3").Value = "'HashMap"
$wsEdit.Range("C3").Value = "'9002"

$wsDelete.Range("A3").Value = "'9809"
$wsDelete.Range("B3").Value = "'HashMap"
$wsDelete.Range("C3").Value = "'9002"

$wsCreate.Range("A3").Value = "'9806"
$wsCreate.Range("B3").Value = "'ListMap"
$wsCreate.Range("C3").Value = "'9005"

# --- Sheet view / selection changes ---
# Update the remembered selection on "Edit" and "Delete" first (this
# temporarily activates each of them), then finish on "Create" so it ends
# up the active tab, matching the target workbook.
$wsEdit.Range("C4").Select()
$wsDelete.Range("C3").Select()
$wsCreate.Range("C3").Select()
